$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the "Obrigatorio" (column E) flag changes from "N" to "S".
# Row 10 is intentionally skipped - it stays "N" per the diff.
$rows = @(2,3,4,5,6,7,8,9,11,12,13,14,15,16,17)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "S"
}
